$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.07172065483994
$ws.Range("D2").Value = 1.07233951787018
$ws.Range("E2").Value = 1.075536613439649
$ws.Range("F2").Value = 1.084898022139389
$ws.Range("I2").Value = 1.056687505477585
$ws.Range("J2").Value = 1.07664382291036
$ws.Range("K2").Value = 1.075034294231135
$ws.Range("L2").Value = 1.078222923148158
$ws.Range("M2").Value = 1.087559850481638
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.072846994044984
$ws.Range("D3").Value = 1.073214745261722
$ws.Range("E3").Value = 1.076516678603262
$ws.Range("F3").Value = 1.08590222190778
$ws.Range("I3").Value = 1.057023721535139
$ws.Range("J3").Value = 1.077427355294653
$ws.Range("K3").Value = 1.075726102761785
$ws.Range("L3").Value = 1.079019922956136
$ws.Range("M3").Value = 1.088382694705964
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.07357612397723
$ws.Range("D4").Value = 1.073781312809023
$ws.Range("E4").Value = 1.077151391187131
$ws.Range("F4").Value = 1.086552567022072
$ws.Range("I4").Value = 1.057240263459714
$ws.Range("J4").Value = 1.077934059826053
$ws.Range("K4").Value = 1.076173343002858
$ws.Range("L4").Value = 1.079535557031641
$ws.Range("M4").Value = 1.088915070758558
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.073882725371217
$ws.Range("D5").Value = 1.07401955402564
$ws.Range("E5").Value = 1.077418354377001
$ws.Range("F5").Value = 1.086826105889157
$ws.Range("I5").Value = 1.057331055177647
$ws.Range("J5").Value = 1.078147007994584
$ws.Range("K5").Value = 1.076361265291357
$ws.Range("L5").Value = 1.079752310608274
$ws.Range("M5").Value = 1.089138867148556
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.07393420950842
$ws.Range("D6").Value = 1.074059559056643
$ws.Range("E6").Value = 1.077463186283321
$ws.Range("F6").Value = 1.086872042100386
$ws.Range("I6").Value = 1.057346285281577
$ws.Range("J6").Value = 1.078182758794213
$ws.Range("K6").Value = 1.076392812544495
$ws.Range("L6").Value = 1.079788703365601
$ws.Range("M6").Value = 1.08917644268711
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.073580220504426
$ws.Range("D7").Value = 1.073784495979531
$ws.Range("E7").Value = 1.077154957851983
$ws.Range("F7").Value = 1.086556221536698
$ws.Range("I7").Value = 1.057241477576161
$ws.Range("J7").Value = 1.077936905526832
$ws.Range("K7").Value = 1.076175854413162
$ws.Range("L7").Value = 1.079538453378516
$ws.Range("M7").Value = 1.088918061194273
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.072101241727776
$ws.Range("D8").Value = 1.072635255715532
$ws.Range("E8").Value = 1.075867717970115
$ws.Range("F8").Value = 1.085237279778731
$ws.Range("I8").Value = 1.056801340791074
$ws.Range("J8").Value = 1.076908681962089
$ws.Range("K8").Value = 1.0752681780887
$ws.Range("L8").Value = 1.078492288958123
$ws.Range("M8").Value = 1.087837946592534
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06949748123197
$ws.Range("D9").Value = 1.070611983483147
$ws.Range("E9").Value = 1.073603628622904
$ws.Range("F9").Value = 1.082917450115906
$ws.Range("I9").Value = 1.05601801666587
$ws.Range("J9").Value = 1.075094582635957
$ws.Range("K9").Value = 1.073665637970376
$ws.Range("L9").Value = 1.076648228150466
$ws.Range("M9").Value = 1.085934207889327
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.067763231388563
$ws.Range("D10").Value = 1.069264396814443
$ws.Range("E10").Value = 1.072097076903652
$ws.Range("F10").Value = 1.081373822127671
$ws.Range("I10").Value = 1.055490599663329
$ws.Range("J10").Value = 1.073883687536448
$ws.Range("K10").Value = 1.072595209072353
$ws.Range("L10").Value = 1.07541847571462
$ws.Range("M10").Value = 1.084664771230001
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.067012654893675
$ws.Range("D11").Value = 1.068681181181303
$ws.Range("E11").Value = 1.071445401361244
$ws.Range("F11").Value = 1.080706112393261
$ws.Range("I11").Value = 1.055260989920944
$ws.Range("J11").Value = 1.073359003086668
$ws.Range("K11").Value = 1.072131213267674
$ws.Range("L11").Value = 1.074885892044215
$ws.Range("M11").Value = 1.084115028323796
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.066733911733052
$ws.Range("D12").Value = 1.068464594071114
$ws.Range("E12").Value = 1.071203441148126
$ws.Range("F12").Value = 1.080458199409261
$ws.Range("I12").Value = 1.055175517178847
$ws.Range("J12").Value = 1.073164058007954
$ws.Range("K12").Value = 1.071958790607647
$ws.Range("L12").Value = 1.074688052937296
$ws.Range("M12").Value = 1.083910819202641
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.066793700634687
$ws.Range("D13").Value = 1.068511050693327
$ws.Range("E13").Value = 1.071255337856291
$ws.Range("F13").Value = 1.080511372868921
$ws.Range("I13").Value = 1.055193859763984
$ws.Range("J13").Value = 1.073205876828632
$ws.Range("K13").Value = 1.071995779198473
$ws.Range("L13").Value = 1.074730490720733
$ws.Range("M13").Value = 1.083954623208136
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.066989612790753
$ws.Range("D14").Value = 1.068663277100383
$ws.Range("E14").Value = 1.071425398790731
$ws.Range("F14").Value = 1.080685617702647
$ws.Range("I14").Value = 1.055253928500516
$ws.Range("J14").Value = 1.073342889965405
$ws.Range("K14").Value = 1.072116962264881
$ws.Range("L14").Value = 1.074869538883131
$ws.Range("M14").Value = 1.084098148539844
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067110327925951
$ws.Range("D15").Value = 1.068757074805724
$ws.Range("E15").Value = 1.071530192349144
$ws.Range("F15").Value = 1.080792989503351
$ws.Range("I15").Value = 1.055290914257538
$ws.Range("J15").Value = 1.073427301137262
$ws.Range("K15").Value = 1.072191617354818
$ws.Range("L15").Value = 1.074955209217598
$ws.Range("M15").Value = 1.084186577896237
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06781305232969
$ws.Range("D16").Value = 1.069303109227889
$ws.Range("E16").Value = 1.072140340627404
$ws.Range("F16").Value = 1.081418150413122
$ws.Range("I16").Value = 1.055505812095308
$ws.Range("J16").Value = 1.073918501554544
$ws.Range("K16").Value = 1.07262599256784
$ws.Range("L16").Value = 1.075453819594523
$ws.Range("M16").Value = 1.084701254427143
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.068253950192574
$ws.Range("D17").Value = 1.069645702287057
$ws.Range("E17").Value = 1.072523250419283
$ws.Range("F17").Value = 1.08181048262537
$ws.Range("I17").Value = 1.055640281294108
$ws.Range("J17").Value = 1.074226522519278
$ws.Range("K17").Value = 1.072898332744543
$ws.Range("L17").Value = 1.075766559960378
$ws.Range("M17").Value = 1.085024079199532
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.068511153718934
$ws.Range("D18").Value = 1.069845559780493
$ws.Range("E18").Value = 1.072746659976167
$ws.Range("F18").Value = 1.082039390056071
$ws.Range("I18").Value = 1.055718595769302
$ws.Range("J18").Value = 1.074406151309374
$ws.Range("K18").Value = 1.073057136584205
$ws.Range("L18").Value = 1.075948967167628
$ws.Range("M18").Value = 1.08521237081192
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.068598859470714
$ws.Range("D19").Value = 1.069913710888865
$ws.Range("E19").Value = 1.072822847837553
$ws.Range("F19").Value = 1.082117452877857
$ws.Range("I19").Value = 1.055745278773268
$ws.Range("J19").Value = 1.07446739421126
$ws.Range("K19").Value = 1.07311127651625
$ws.Range("L19").Value = 1.076011161788709
$ws.Range("M19").Value = 1.085276572259018
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.068206642380542
$ws.Range("D20").Value = 1.069608942317018
$ws.Range("E20").Value = 1.072482161122809
$ws.Range("F20").Value = 1.081768382191987
$ws.Range("I20").Value = 1.055625866343882
$ws.Range("J20").Value = 1.074193478366382
$ws.Range("K20").Value = 1.072869118141283
$ws.Range("L20").Value = 1.075733006804408
$ws.Range("M20").Value = 1.084989443839705
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066931920046812
$ws.Range("D21").Value = 1.068618448968811
$ws.Range("E21").Value = 1.071375317291386
$ws.Range("F21").Value = 1.080634304065065
$ws.Range("I21").Value = 1.055236244874066
$ws.Range("J21").Value = 1.073302544499777
$ws.Range("K21").Value = 1.072081278912934
$ws.Range("L21").Value = 1.074828593045499
$ws.Range("M21").Value = 1.084055884185259
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.066130765846991
$ws.Range("D22").Value = 1.067995947839466
$ws.Range("E22").Value = 1.070679985695508
$ws.Range("F22").Value = 1.079921866799209
$ws.Range("I22").Value = 1.054990201396907
$ws.Range("J22").Value = 1.072742067100718
$ws.Range("K22").Value = 1.071585505276173
$ws.Range("L22").Value = 1.074259872329377
$ws.Range("M22").Value = 1.083468859656524
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.066555443213756
$ws.Range("D23").Value = 1.068325922582863
$ws.Range("E23").Value = 1.071048538592219
$ws.Range("F23").Value = 1.080299486057189
$ws.Range("I23").Value = 1.055120735393103
$ws.Range("J23").Value = 1.073039216338018
$ws.Range("K23").Value = 1.071848364777689
$ws.Range("L23").Value = 1.074561369474736
$ws.Range("M23").Value = 1.083780057960426
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.068228018635274
$ws.Range("D24").Value = 1.069625552476247
$ws.Range("E24").Value = 1.072500727407059
$ws.Range("F24").Value = 1.081787405359504
$ws.Range("I24").Value = 1.055632380207449
$ws.Range("J24").Value = 1.074208409703585
$ws.Range("K24").Value = 1.07288231910989
$ws.Range("L24").Value = 1.075748168058895
$ws.Range("M24").Value = 1.085005094088302
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.070170334259381
$ws.Range("D25").Value = 1.071134827494607
$ws.Range("E25").Value = 1.074188450890802
$ws.Range("F25").Value = 1.083516667915172
$ws.Range("I25").Value = 1.056221441994903
$ws.Range("J25").Value = 1.075563835077407
$ws.Range("K25").Value = 1.074080298709716
$ws.Range("L25").Value = 1.077125029904549
$ws.Range("M25").Value = 1.086426420179063
